# Update "want-to-go" counts and the cover image URL for the
# "丽水·龙泉ACG动漫游戏博览会" event on both the "展览" and "全部类型" sheets.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    # Row 3: 丽水·龙泉ACG动漫游戏博览会
    $ws.Range("F3").Value = 85
    $ws.Range("I3").Value = "//i1.hdslb.com/bfs/openplatform/202406/yDEGBE9v1717471354651.png"

    # Row 4: 丽水·CCAC动漫七夕（回馈展）
    $ws.Range("F4").Value = 14

    # Row 5: 丽水·AEO纯白礼赞动漫嘉年华
    $ws.Range("F5").Value = 26
}
